# Update "想去人数" (want-to-go count) values in column F on both the
# "展览" sheet and the "全部类型" sheet, as per the upstream data refresh.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (row -> new F value)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value  = 2828
$ws1.Range("F8").Value  = 1680
$ws1.Range("F9").Value  = 1822
$ws1.Range("F12").Value = 747
$ws1.Range("F20").Value = 6637
$ws1.Range("F22").Value = 1536
$ws1.Range("F28").Value = 68
$ws1.Range("F34").Value = 778
$ws1.Range("F35").Value = 1470
$ws1.Range("F37").Value = 142
$ws1.Range("F42").Value = 161

# Sheet "全部类型" (row -> new F value)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 2828
$ws4.Range("F11").Value = 1680
$ws4.Range("F12").Value = 1822
$ws4.Range("F15").Value = 747
$ws4.Range("F23").Value = 6637
$ws4.Range("F25").Value = 1536
$ws4.Range("F32").Value = 68
$ws4.Range("F38").Value = 778
$ws4.Range("F39").Value = 1470
$ws4.Range("F41").Value = 142
$ws4.Range("F49").Value = 161
